# Fix typo in predictor labels: "per capita" -> "per cap."
# Also normalize the bracket typo on the Livestock AB Consumption row
# ( "[kg per capita)" -> "(kg per cap.)" ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "ln(GDP [dollars per capita])"               = "ln(GDP [dollars per cap.])"
    "ln(Migrant Population [per capita])"        = "ln(Migrant Population [per cap.])"
    "ln(Tourism - Inbound [per capita])"         = "ln(Tourism - Inbound [per cap.])"
    "ln(ProMed Mentions [per capita])"           = "ln(ProMed Mentions [per cap.])"
    "ln(Publication Bias Index [per capita])"    = "ln(Publication Bias Index [per cap.])"
    "ln(AB Exports [dollars per capita])"        = "ln(AB Exports [dollars per cap.])"
    "Livestock AB Consumption [kg per capita)"   = "Livestock AB Consumption (kg per cap.)"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
